$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use an unused cell (outside the populated data columns, but inside the
# sheet's existing dimension so no structural column gets inserted) as a
# scratch cell to funnel values through as genuine Text, then
# PasteSpecial-values them into the target cells. This avoids Excel's
# automatic "looks like a number" type coercion that a plain
# Range.Value = "34.32" assignment would trigger (which would turn the
# shared-string text cell into a numeric cell and touch its style).
$scratch = $ws.Range("E1")

# Each pair is: target cell address, new text value (old value, per the
# commit diff, is shown in the comment for reference).
$pairs = @(
    @("B11", "34.32"),   # 34.3  -> 34.32
    @("C11", "7.02"),    # 7     -> 7.02
    @("D11", "41.34"),   # 41.3  -> 41.34
    @("B12", "14.44"),   # 14.4  -> 14.44
    @("C12", "26.57"),   # 26.6  -> 26.57
    @("D12", "41.01"),   # 41    -> 41.01
    @("B33", "24.19"),   # 24.2  -> 24.19
    @("C33", "2.75"),    # 2.7   -> 2.75
    @("D33", "26.94"),   # 26.9  -> 26.94
    @("B34", "17.15"),   # 17.2  -> 17.15
    @("C34", "35.95"),   # 35.9  -> 35.95
    @("B36", "89.49"),   # 89.5  -> 89.49
    @("C36", "10.17"),   # 10.2  -> 10.17
    @("D36", "99.65"),   # 99.7  -> 99.65
    @("B40", "18.66"),   # 18.7  -> 18.66
    @("C40", "32.39"),   # 32.4  -> 32.39
    @("D40", "51.05")    # 51.1  -> 51.05
)

foreach ($pair in $pairs) {
    $addr = $pair[0]
    $newValue = $pair[1]

    $scratch.NumberFormat = "@"
    $scratch.Value = $newValue
    $scratch.Copy()
    $ws.Range($addr).PasteSpecial(-4163)
    $scratch.Clear()
}
